$wb = $excel.ActiveWorkbook

# Sheet "OFF" - row 2 (H) updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 332
$wsOff.Range("C2").Value = 215
$wsOff.Range("D2").Value = 147
$wsOff.Range("E2").Value = 64

# Sheet "DEF" - row 2 (H) updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 452
$wsDef.Range("C2").Value = 341
$wsDef.Range("D2").Value = 88
$wsDef.Range("E2").Value = 54
